$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2022-05-08) updated values
$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 5.586269137925634

# Row 3 (2022-04-01) updated values
$ws.Range("B3").Value = 0.01293466051926884
$ws.Range("C3").Value = 0.306821227259698
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 1.566735691613601
